$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.046.19'
$ws.Range('E2').Value = '  +2.81%  '
$ws.Range('D3').Value = '1.652.70'
$ws.Range('E3').Value = '  +3.60%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.29'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.63%  '
$ws.Range('E6').Value = '  +1.61%  '
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('E8').Value = '  +1.76%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0615'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.68%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.95'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.79%  '
$ws.Range('E11').Value = '  +1.37%  '
$ws.Range('D12').Value = '1.886.51'
$ws.Range('E12').Value = '  +3.59%  '
$ws.Range('D13').Value = '1.647.62'
$ws.Range('E13').Value = '  +3.29%  '
$ws.Range('E14').Value = '  +2.25%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.519'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.02%  '
$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.35'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.04%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '240.37'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +4.58%  '
$ws.Range('D18').Value = '27.032.69'
$ws.Range('E18').Value = '  +2.75%  '
$ws.Range('E19').Value = '  +2.82%  '
$ws.Range('E20').Value = '  +1.20%  '
$ws.Range('E21').Value = '  -0.04%  '
$ws.Range('E22').Value = '  +4.36%  '
$ws.Range('E23').Value = '  +2.70%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.26'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.80%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.05'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.37%  '
$ws.Range('E26').Value = '  -0.18%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.12'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.29%  '
$ws.Range('E28').Value = '  +1.96%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.84'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.17%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0498'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.91%  '
$ws.Range('E31').Value = '  +2.01%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.31'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.23%  '
$ws.Range('D33').Value = '1.521.44'
$ws.Range('E33').Value = '  +1.41%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.09'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.57%  '
$ws.Range('E35').Value = '  +8.57%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.41'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.14%  '
$ws.Range('E37').Value = '  +2.57%  '
$ws.Range('B38').Value = 'ARBITRUM'
$ws.Range('C38').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.888'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +8.92%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0169'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.31%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.96'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.99%  '
$ws.Range('E41').Value = '  -0.13%  '
$ws.Range('E42').Value = '  +4.02%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '65.75'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +8.53%  '
$ws.Range('D44').Value = '1.793.04'
$ws.Range('E44').Value = '  +3.38%  '
$ws.Range('E45').Value = '  +2.07%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.916'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.26%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '89.58'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.28%  '
$ws.Range('D48').Value = '0.0₆0105'
$ws.Range('E48').Value = '  +1.02%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.53'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.97%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0508'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.45%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0977'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.18%  '
